$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells are treated as text so values like "1.00" or "521.83" are preserved exactly
$ws.Range("D2:E51").NumberFormat = "@"

# Update price (D) and volume/1h (E) columns for rows 2-51 per latest crypto data refresh
$ws.Range("D2").Value = "58.176.79"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "2.478.23"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").Value = "521.83"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").Value = "132.24"
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "0.557"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "2.510.56"
$ws.Range("E9").Value = "  +1.56%  "
$ws.Range("D10").Value = "0.0972"
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("D12").Value = "5.14"
$ws.Range("E12").Value = "  -2.65%  "
$ws.Range("D13").Value = "0.331"
$ws.Range("E13").Value = "  -2.19%  "
$ws.Range("D14").Value = "2.957.61"
$ws.Range("E14").Value = "  +1.88%  "
$ws.Range("D15").Value = "58.232.92"
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("D16").Value = "22.02"
$ws.Range("E16").Value = "  -1.16%  "
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").Value = "2.501.45"
$ws.Range("E18").Value = "  +1.46%  "
$ws.Range("D19").Value = "10.59"
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("D20").Value = "320.88"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("D22").Value = "6.11"
$ws.Range("E22").Value = "  +7.01%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "64.27"
$ws.Range("E24").Value = "  -0.50%  "
$ws.Range("D25").Value = "0.402"
$ws.Range("E25").Value = "  -1.35%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  +0.30%  "
$ws.Range("D28").Value = "7.37"
$ws.Range("E28").Value = "  +0.67%  "
$ws.Range("D29").Value = "0.0₃0751"
$ws.Range("E29").Value = "  +0.59%  "
$ws.Range("D30").Value = "1.71"
$ws.Range("E30").Value = "  +1.55%  "
$ws.Range("D33").Value = "6.26"
$ws.Range("E33").Value = "  +0.62%  "
$ws.Range("E34").Value = "  -0.01%  "
$ws.Range("D35").Value = "0.993"
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("D36").Value = "18.04"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  -8.74%  "
$ws.Range("D38").Value = "3.94"
$ws.Range("E38").Value = "  -0.55%  "
$ws.Range("E39").Value = "  +0.20%  "
$ws.Range("D40").Value = "36.08"
$ws.Range("E40").Value = "  -0.60%  "
$ws.Range("D41").Value = "0.771"
$ws.Range("E41").Value = "  -3.02%  "
$ws.Range("D42").Value = "278.30"
$ws.Range("E42").Value = "  +1.61%  "
$ws.Range("E43").Value = "  +1.07%  "
$ws.Range("D44").Value = "4.99"
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("D45").Value = "0.594"
$ws.Range("E45").Value = "  +0.93%  "
$ws.Range("D46").Value = "123.32"
$ws.Range("E46").Value = "  -0.23%  "
$ws.Range("D47").Value = "0.0918"
$ws.Range("E47").Value = "  +1.02%  "
$ws.Range("D48").Value = "0.0499"
$ws.Range("E48").Value = "  +2.47%  "
$ws.Range("D49").Value = "17.62"
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("E50").Value = "  +0.41%  "
$ws.Range("D51").Value = "16.75"
$ws.Range("E51").Value = "  -0.52%  "

# Row 31/32 swap: Monero now ranks above Fetch.AI, update coin name, link, price, and volume
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "167.22"
$ws.Range("E31").Value = "  -0.21%  "

$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "1.19"
$ws.Range("E32").Value = "  +2.09%  "
